$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A14").Value = "bal4"

$ws.Range("C14").Value = 10
$ws.Range("D14").Value = 0.02

$ws.Range("C15").Value = 11
$ws.Range("D15").Value = 0.12

$ws.Range("C16").Value = 12
$ws.Range("D16").Value = 0.56000000000000005

$ws.Range("C17").Value = 13
$ws.Range("D17").Value = 2.52

$ws.Range("C18").Value = 14
$ws.Range("D18").Value = 11.39

$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 54.75

$ws.Range("D20").Select()
